$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: C114 should reference "isolated" not "isoalted"
$ws.Range("C114").Value = "isolated"

# Row 118
$ws.Range("A118").Value = 5
$ws.Range("B118").Value = ",0"
$ws.Range("C118").Value = "social"
$ws.Range("D118").Value = 1
$ws.Range("E118").Value = 1
$ws.Range("F118").Value = 16
$ws.Range("G118").Value = 6
$ws.Range("H118").Formula = "=F118-G118"
$ws.Range("I118").Value = 3
$ws.Range("J118").Value = 7
$ws.Range("K118").Value = 5
$ws.Range("L118").Value = 7

# Row 119
$ws.Range("A119").Value = 5
$ws.Range("B119").Value = ",0"
$ws.Range("C119").Value = "social"
$ws.Range("D119").Value = 2
$ws.Range("E119").Value = 3
$ws.Range("F119").Value = 3
$ws.Range("G119").Value = 0
$ws.Range("H119").Formula = "=F119-G119"
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 1
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 3

# Row 120
$ws.Range("A120").Value = 5
$ws.Range("B120").Value = ",0"
$ws.Range("C120").Value = "social"
$ws.Range("D120").Value = "both"
$ws.Range("E120").Formula = "=SUM(E118:E119)"
$ws.Range("F120").Formula = "=SUM(F118:F119)"
$ws.Range("G120").Formula = "=SUM(G118:G119)"
$ws.Range("H120").Formula = "=SUM(H118:H119)"
$ws.Range("I120").Formula = "=SUM(I118:I119)"
$ws.Range("J120").Formula = "=SUM(J118:J119)"
$ws.Range("K120").Formula = "=SUM(K118:K119)"
$ws.Range("L120").Formula = "=SUM(L118:L119)"

# Row 121
$ws.Range("A121").Value = 5
$ws.Range("B121").Value = ",1"
$ws.Range("C121").Value = "social"
$ws.Range("D121").Value = 1
$ws.Range("E121").Value = 4
$ws.Range("F121").Value = 6
$ws.Range("G121").Value = 0
$ws.Range("H121").Formula = "=F121-G121"
$ws.Range("I121").Value = 1
$ws.Range("J121").Value = 3
$ws.Range("K121").Value = 1
$ws.Range("L121").Value = 4

# Row 122
$ws.Range("A122").Value = 5
$ws.Range("B122").Value = ",1"
$ws.Range("C122").Value = "social"
$ws.Range("D122").Value = 2
$ws.Range("E122").Value = 2
$ws.Range("F122").Value = 8
$ws.Range("G122").Value = 3
$ws.Range("H122").Formula = "=F122-G122"
$ws.Range("I122").Value = 1
$ws.Range("J122").Value = 2
$ws.Range("K122").Value = 1
$ws.Range("L122").Value = 3

# Row 123
$ws.Range("A123").Value = 5
$ws.Range("B123").Value = ",1"
$ws.Range("C123").Value = "social"
$ws.Range("D123").Value = "both"
$ws.Range("E123").Formula = "=SUM(E121:E122)"
$ws.Range("F123").Formula = "=SUM(F121:F122)"
$ws.Range("G123").Formula = "=SUM(G121:G122)"
$ws.Range("H123").Formula = "=SUM(H121:H122)"
$ws.Range("I123").Formula = "=SUM(I121:I122)"
$ws.Range("J123").Formula = "=SUM(J121:J122)"
$ws.Range("K123").Formula = "=SUM(K121:K122)"
$ws.Range("L123").Formula = "=SUM(L121:L122)"

# Row 124
$ws.Range("A124").Value = 5
$ws.Range("B124").Value = ",3"
$ws.Range("C124").Value = "isolated"
$ws.Range("D124").Value = 1
$ws.Range("E124").Value = 4
$ws.Range("F124").Value = 14
$ws.Range("G124").Value = 1
$ws.Range("H124").Formula = "=F124-G124"
$ws.Range("I124").Value = 1
$ws.Range("J124").Value = 2
$ws.Range("K124").Value = 11
$ws.Range("L124").Value = 15

# Row 125
$ws.Range("A125").Value = 5
$ws.Range("B125").Value = ",3"
$ws.Range("C125").Value = "isolated"
$ws.Range("D125").Value = 2
$ws.Range("E125").Value = 0
$ws.Range("F125").Value = 2
$ws.Range("G125").Value = 1
$ws.Range("H125").Formula = "=F125-G125"
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 1
$ws.Range("L125").Value = 1

# Row 126
$ws.Range("A126").Value = 5
$ws.Range("B126").Value = ",3"
$ws.Range("C126").Value = "isolated"
$ws.Range("D126").Value = "both"
$ws.Range("E126").Formula = "=SUM(E124:E125)"
$ws.Range("F126").Formula = "=SUM(F124:F125)"
$ws.Range("G126").Formula = "=SUM(G124:G125)"
$ws.Range("H126").Formula = "=SUM(H124:H125)"
$ws.Range("I126").Formula = "=SUM(I124:I125)"
$ws.Range("J126").Formula = "=SUM(J124:J125)"
$ws.Range("K126").Formula = "=SUM(K124:K125)"
$ws.Range("L126").Formula = "=SUM(L124:L125)"

# Row 127
$ws.Range("A127").Value = 5
$ws.Range("B127").Value = ",4"
$ws.Range("C127").Value = "social"
$ws.Range("D127").Value = 1
$ws.Range("E127").Value = 1
$ws.Range("F127").Value = 9
$ws.Range("G127").Value = 7
$ws.Range("H127").Formula = "=F127-G127"
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 1
$ws.Range("K127").Value = 1
$ws.Range("L127").Value = 2

# Row 128
$ws.Range("A128").Value = 5
$ws.Range("B128").Value = ",4"
$ws.Range("C128").Value = "social"
$ws.Range("D128").Value = 2
$ws.Range("E128").Value = 1
$ws.Range("F128").Value = 2
$ws.Range("G128").Value = 0
$ws.Range("H128").Formula = "=F128-G128"
$ws.Range("I128").Value = 1
$ws.Range("J128").Value = 1
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 1

# Row 129
$ws.Range("A129").Value = 5
$ws.Range("B129").Value = ",4"
$ws.Range("C129").Value = "social"
$ws.Range("D129").Value = "both"
$ws.Range("E129").Formula = "=SUM(E127:E128)"
$ws.Range("F129").Formula = "=SUM(F127:F128)"
$ws.Range("G129").Formula = "=SUM(G127:G128)"
$ws.Range("H129").Formula = "=SUM(H127:H128)"
$ws.Range("I129").Formula = "=SUM(I127:I128)"
$ws.Range("J129").Formula = "=SUM(J127:J128)"
$ws.Range("K129").Formula = "=SUM(K127:K128)"
$ws.Range("L129").Formula = "=SUM(L127:L128)"

# Row 130
$ws.Range("A130").Value = 5
$ws.Range("B130").Value = ",5"
$ws.Range("C130").Value = "social"
$ws.Range("D130").Value = 1
$ws.Range("E130").Value = 1
$ws.Range("F130").Value = 10
$ws.Range("G130").Value = 8
$ws.Range("H130").Formula = "=F130-G130"
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 1
$ws.Range("K130").Value = 2
$ws.Range("L130").Value = 3

# Row 131
$ws.Range("A131").Value = 5
$ws.Range("B131").Value = ",5"
$ws.Range("C131").Value = "social"
$ws.Range("D131").Value = 2
$ws.Range("E131").Value = 0
$ws.Range("F131").Value = 3
$ws.Range("G131").Value = 2
$ws.Range("H131").Formula = "=F131-G131"
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 1
$ws.Range("K131").Value = 1
$ws.Range("L131").Value = 1

# Row 132
$ws.Range("A132").Value = 5
$ws.Range("B132").Value = ",5"
$ws.Range("C132").Value = "social"
$ws.Range("D132").Value = "both"
$ws.Range("E132").Formula = "=SUM(E130:E131)"
$ws.Range("F132").Formula = "=SUM(F130:F131)"
$ws.Range("G132").Formula = "=SUM(G130:G131)"
$ws.Range("H132").Formula = "=SUM(H130:H131)"
$ws.Range("I132").Formula = "=SUM(I130:I131)"
$ws.Range("J132").Formula = "=SUM(J130:J131)"
$ws.Range("K132").Formula = "=SUM(K130:K131)"
$ws.Range("L132").Formula = "=SUM(L130:L131)"

# Row 133
$ws.Range("A133").Value = 5
$ws.Range("B133").Value = ",6"
$ws.Range("C133").Value = "isolated"
$ws.Range("D133").Value = 1
$ws.Range("E133").Value = 4
$ws.Range("F133").Value = 12
$ws.Range("G133").Value = 3
$ws.Range("H133").Formula = "=F133-G133"
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 5
$ws.Range("K133").Value = 5
$ws.Range("L133").Value = 9

# Row 134
$ws.Range("A134").Value = 5
$ws.Range("B134").Value = ",6"
$ws.Range("C134").Value = "isolated"
$ws.Range("D134").Value = 2
$ws.Range("E134").Value = 2
$ws.Range("F134").Value = 12
$ws.Range("G134").Value = 3
$ws.Range("H134").Formula = "=F134-G134"
$ws.Range("I134").Value = 3
$ws.Range("J134").Value = 4
$ws.Range("K134").Value = 5
$ws.Range("L134").Value = 8

# Row 135
$ws.Range("A135").Value = 5
$ws.Range("B135").Value = ",6"
$ws.Range("C135").Value = "isolated"
$ws.Range("D135").Value = "both"
$ws.Range("E135").Formula = "=SUM(E133:E134)"
$ws.Range("F135").Formula = "=SUM(F133:F134)"
$ws.Range("G135").Formula = "=SUM(G133:G134)"
$ws.Range("H135").Formula = "=SUM(H133:H134)"
$ws.Range("I135").Formula = "=SUM(I133:I134)"
$ws.Range("J135").Formula = "=SUM(J133:J134)"
$ws.Range("K135").Formula = "=SUM(K133:K134)"
$ws.Range("L135").Formula = "=SUM(L133:L134)"

# Row 136
$ws.Range("A136").Value = 5
$ws.Range("B136").Value = ",7"
$ws.Range("C136").Value = "social"
$ws.Range("D136").Value = 1
$ws.Range("E136").Value = 4
$ws.Range("F136").Value = 27
$ws.Range("G136").Value = 17
$ws.Range("H136").Formula = "=F136-G136"
$ws.Range("I136").Value = 4
$ws.Range("J136").Value = 6
$ws.Range("K136").Value = 2
$ws.Range("L136").Value = 5

# Row 137
$ws.Range("A137").Value = 5
$ws.Range("B137").Value = ",7"
$ws.Range("C137").Value = "social"
$ws.Range("D137").Value = 2
$ws.Range("E137").Value = 1
$ws.Range("F137").Value = 3
$ws.Range("G137").Value = 2
$ws.Range("H137").Formula = "=F137-G137"
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 1
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 1

# Row 138
$ws.Range("A138").Value = 5
$ws.Range("B138").Value = ",7"
$ws.Range("C138").Value = "social"
$ws.Range("D138").Value = "both"
$ws.Range("E138").Formula = "=SUM(E136:E137)"
$ws.Range("F138").Formula = "=SUM(F136:F137)"
$ws.Range("G138").Formula = "=SUM(G136:G137)"
$ws.Range("H138").Formula = "=SUM(H136:H137)"
$ws.Range("I138").Formula = "=SUM(I136:I137)"
$ws.Range("J138").Formula = "=SUM(J136:J137)"
$ws.Range("K138").Formula = "=SUM(K136:K137)"
$ws.Range("L138").Formula = "=SUM(L136:L137)"

# Row 139
$ws.Range("A139").Value = 5
$ws.Range("B139").Value = ",8"
$ws.Range("C139").Value = "isolated"
$ws.Range("D139").Value = 1
$ws.Range("E139").Value = 3
$ws.Range("F139").Value = 18
$ws.Range("G139").Value = 12
$ws.Range("H139").Formula = "=F139-G139"
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 2
$ws.Range("L139").Value = 5

# Row 140
$ws.Range("A140").Value = 5
$ws.Range("B140").Value = ",8"
$ws.Range("C140").Value = "isolated"
$ws.Range("D140").Value = 2
$ws.Range("E140").Value = 0
$ws.Range("F140").Value = 4
$ws.Range("G140").Value = 1
$ws.Range("H140").Formula = "=F140-G140"
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 2
$ws.Range("K140").Value = 3
$ws.Range("L140").Value = 3

# Row 141
$ws.Range("A141").Value = 5
$ws.Range("B141").Value = ",8"
$ws.Range("C141").Value = "isolated"
$ws.Range("D141").Value = "both"
$ws.Range("E141").Formula = "=SUM(E139:E140)"
$ws.Range("F141").Formula = "=SUM(F139:F140)"
$ws.Range("G141").Formula = "=SUM(G139:G140)"
$ws.Range("H141").Formula = "=SUM(H139:H140)"
$ws.Range("I141").Formula = "=SUM(I139:I140)"
$ws.Range("J141").Formula = "=SUM(J139:J140)"
$ws.Range("K141").Formula = "=SUM(K139:K140)"
$ws.Range("L141").Formula = "=SUM(L139:L140)"

# Row 142
$ws.Range("A142").Value = 5
$ws.Range("B142").Value = ",9"
$ws.Range("C142").Value = "social"
$ws.Range("D142").Value = 1
$ws.Range("E142").Value = 1
$ws.Range("F142").Value = 15
$ws.Range("G142").Value = 7
$ws.Range("H142").Formula = "=F142-G142"
$ws.Range("I142").Value = 5
$ws.Range("J142").Value = 6
$ws.Range("K142").Value = 2
$ws.Range("L142").Value = 3

# Row 143
$ws.Range("A143").Value = 5
$ws.Range("B143").Value = ",9"
$ws.Range("C143").Value = "social"
$ws.Range("D143").Value = 2
$ws.Range("E143").Value = 1
$ws.Range("F143").Value = 6
$ws.Range("G143").Value = 3
$ws.Range("H143").Formula = "=F143-G143"
$ws.Range("I143").Value = 1
$ws.Range("J143").Value = 1
$ws.Range("K143").Value = 1
$ws.Range("L143").Value = 2

# Row 144
$ws.Range("A144").Value = 5
$ws.Range("B144").Value = ",9"
$ws.Range("C144").Value = "social"
$ws.Range("D144").Value = "both"
$ws.Range("E144").Formula = "=SUM(E142:E143)"
$ws.Range("F144").Formula = "=SUM(F142:F143)"
$ws.Range("G144").Formula = "=SUM(G142:G143)"
$ws.Range("H144").Formula = "=SUM(H142:H143)"
$ws.Range("I144").Formula = "=SUM(I142:I143)"
$ws.Range("J144").Formula = "=SUM(J142:J143)"
$ws.Range("K144").Formula = "=SUM(K142:K143)"
$ws.Range("L144").Formula = "=SUM(L142:L143)"

# Row 145
$ws.Range("A145").Value = 5
$ws.Range("B145").Value = ":0"
$ws.Range("C145").Value = "isolated"
$ws.Range("D145").Value = 1
$ws.Range("E145").Value = 3
$ws.Range("F145").Value = 31
$ws.Range("G145").Value = 18
$ws.Range("H145").Formula = "=F145-G145"
$ws.Range("I145").Value = 3
$ws.Range("J145").Value = 7
$ws.Range("K145").Value = 7
$ws.Range("L145").Value = 10

# Row 146
$ws.Range("A146").Value = 5
$ws.Range("B146").Value = ":0"
$ws.Range("C146").Value = "isolated"
$ws.Range("D146").Value = 2
$ws.Range("E146").Value = 0
$ws.Range("F146").Value = 3
$ws.Range("G146").Value = 2
$ws.Range("H146").Formula = "=F146-G146"
$ws.Range("I146").Value = 0
$ws.Range("J146").Value = 1
$ws.Range("K146").Value = 1
$ws.Range("L146").Value = 1

# Row 147
$ws.Range("A147").Value = 5
$ws.Range("B147").Value = ":0"
$ws.Range("C147").Value = "isolated"
$ws.Range("D147").Value = "both"
$ws.Range("E147").Formula = "=SUM(E145:E146)"
$ws.Range("F147").Formula = "=SUM(F145:F146)"
$ws.Range("G147").Formula = "=SUM(G145:G146)"
$ws.Range("H147").Formula = "=SUM(H145:H146)"
$ws.Range("I147").Formula = "=SUM(I145:I146)"
$ws.Range("J147").Formula = "=SUM(J145:J146)"
$ws.Range("K147").Formula = "=SUM(K145:K146)"
$ws.Range("L147").Formula = "=SUM(L145:L146)"

# Row 148
$ws.Range("A148").Value = 5
$ws.Range("B148").Value = ":1"
$ws.Range("C148").Value = "isolated"
$ws.Range("D148").Value = 1
$ws.Range("E148").Value = 2
$ws.Range("F148").Value = 20
$ws.Range("G148").Value = 13
$ws.Range("H148").Formula = "=F148-G148"
$ws.Range("I148").Value = 4
$ws.Range("J148").Value = 4
$ws.Range("K148").Value = 1
$ws.Range("L148").Value = 3

# Row 149
$ws.Range("A149").Value = 5
$ws.Range("B149").Value = ":1"
$ws.Range("C149").Value = "isolated"
$ws.Range("D149").Value = 2
$ws.Range("E149").Value = 1
$ws.Range("F149").Value = 1
$ws.Range("G149").Value = 0
$ws.Range("H149").Formula = "=F149-G149"
$ws.Range("I149").Value = 0
$ws.Range("J149").Value = 1
$ws.Range("K149").Value = 0
$ws.Range("L149").Value = 1

# Row 150
$ws.Range("A150").Value = 5
$ws.Range("B150").Value = ":1"
$ws.Range("C150").Value = "isolated"
$ws.Range("D150").Value = "both"
$ws.Range("E150").Formula = "=SUM(E148:E149)"
$ws.Range("F150").Formula = "=SUM(F148:F149)"
$ws.Range("G150").Formula = "=SUM(G148:G149)"
$ws.Range("H150").Formula = "=SUM(H148:H149)"
$ws.Range("I150").Formula = "=SUM(I148:I149)"
$ws.Range("J150").Formula = "=SUM(J148:J149)"
$ws.Range("K150").Formula = "=SUM(K148:K149)"
$ws.Range("L150").Formula = "=SUM(L148:L149)"

# Row 151
$ws.Range("A151").Value = 5
$ws.Range("B151").Value = ":7"
$ws.Range("C151").Value = "social"
$ws.Range("D151").Value = 1
$ws.Range("E151").Value = 2
$ws.Range("F151").Value = 5
$ws.Range("G151").Value = 3
$ws.Range("H151").Formula = "=F151-G151"
$ws.Range("I151").Value = 0
$ws.Range("J151").Value = 0
$ws.Range("K151").Value = 0
$ws.Range("L151").Value = 2

# Row 152
$ws.Range("A152").Value = 5
$ws.Range("B152").Value = ":7"
$ws.Range("C152").Value = "social"
$ws.Range("D152").Value = 2
$ws.Range("E152").Value = 0
$ws.Range("F152").Value = 7
$ws.Range("G152").Value = 4
$ws.Range("H152").Formula = "=F152-G152"
$ws.Range("I152").Value = 0
$ws.Range("J152").Value = 1
$ws.Range("K152").Value = 3
$ws.Range("L152").Value = 3

# Row 153
$ws.Range("A153").Value = 5
$ws.Range("B153").Value = ":7"
$ws.Range("C153").Value = "social"
$ws.Range("D153").Value = "both"
$ws.Range("E153").Formula = "=SUM(E151:E152)"
$ws.Range("F153").Formula = "=SUM(F151:F152)"
$ws.Range("G153").Formula = "=SUM(G151:G152)"
$ws.Range("H153").Formula = "=SUM(H151:H152)"
$ws.Range("I153").Formula = "=SUM(I151:I152)"
$ws.Range("J153").Formula = "=SUM(J151:J152)"
$ws.Range("K153").Formula = "=SUM(K151:K152)"
$ws.Range("L153").Formula = "=SUM(L151:L152)"

# Update view selection to match target state
$ws.Range("K160").Select()
